$wb = $excel.ActiveWorkbook

# Kunden: remove the title row (row 1); header row becomes row 1
$ws = $wb.Worksheets.Item("Kunden")
$ws.Rows("1").Delete()
$ws.Range("D16").Select()

# Produkte: remove the title row and header row (rows 1-2)
$ws = $wb.Worksheets.Item("Produkte")
$ws.Rows("1:2").Delete()
$ws.Range("E18:E19").Select()

# Bestellungen: remove the title row and header row (rows 1-2)
$ws = $wb.Worksheets.Item("Bestellungen")
$ws.Rows("1:2").Delete()
$ws.Range("K20").Select()

# Bestellpositionen: remove the title row and header row (rows 1-2)
$ws = $wb.Worksheets.Item("Bestellpositionen")
$ws.Rows("1:2").Delete()
$ws.Range("L20:L21").Select()

# Städte: remove the title row and header row (rows 1-2)
$ws = $wb.Worksheets.Item("Städte")
$ws.Rows("1:2").Delete()
$ws.Range("E12").Select()

# Hersteller: remove the title row and header row (rows 1-2); this becomes the active sheet
$ws = $wb.Worksheets.Item("Hersteller")
$ws.Rows("1:2").Delete()
$ws.Activate()
$ws.Range("I46").Select()
